$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2399.7917
$ws.Range("I15").Value = 2399.7917
$ws.Range("K15").Value = 7199.375100000001
$ws.Range("M15").Value = -7030.375100000001
$ws.Range("H32").Value = 1927.6666
$ws.Range("I32").Value = 1914.4445
$ws.Range("K32").Value = 1914.4445
$ws.Range("M32").Value = -1588.4445
$ws.Range("H40").Value = 4974.3887
$ws.Range("I40").Value = 1715.6
$ws.Range("K40").Value = 1715.6
$ws.Range("M40").Value = -1540.6
$ws.Range("H70").Value = 4636.364
$ws.Range("J70").Value = 5250
$ws.Range("L70").Value = 15750
$ws.Range("N70").Value = -16290
$ws.Range("H73").Value = 4636.364
$ws.Range("J73").Value = 5250
$ws.Range("L73").Value = 15750
$ws.Range("N73").Value = -17622
$ws.Range("H99").Value = 536.75
$ws.Range("J99").Value = 530
$ws.Range("L99").Value = 1590
$ws.Range("N99").Value = -4586
$ws.Range("H111").Value = 674.875
$ws.Range("I111").Value = 566.6667
$ws.Range("J111").Value = 999.5
$ws.Range("K111").Value = 1700.0001
$ws.Range("L111").Value = 2998.5
$ws.Range("M111").Value = 1366.9999
$ws.Range("N111").Value = -9132.5
$ws.Range("H116").Value = 14523.375
$ws.Range("I116").Value = 19705
$ws.Range("J116").Value = 11414.4
$ws.Range("K116").Value = 19705
$ws.Range("L116").Value = 11414.4
$ws.Range("M116").Value = -16263
$ws.Range("N116").Value = -18298.4
$ws.Range("H132").Value = 34499.668
$ws.Range("I132").Value = 34499.668
$ws.Range("K132").Value = 103499.004
$ws.Range("M132").Value = -100969.004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1359076.4
$ws.Range("I32").Value = 1544842.1
$ws.Range("K32").Value = 1544842.1
$ws.Range("M32").Value = -1544555.1
$ws.Range("H61").Value = 6669899
$ws.Range("I61").Value = 2610.2222
$ws.Range("K61").Value = 2610.2222
$ws.Range("M61").Value = -2398.2222
$ws.Range("H136").Value = 6669899
$ws.Range("I136").Value = 2610.2222
$ws.Range("K136").Value = 7830.6666
$ws.Range("M136").Value = -5280.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 411.16666
$ws.Range("I22").Value = 230.14285
$ws.Range("K22").Value = 230.14285
$ws.Range("M22").Value = -57.14285000000001
$ws.Range("H86").Value = 2741.4
$ws.Range("I86").Value = 1638.091
$ws.Range("J86").Value = 4089.889
$ws.Range("K86").Value = 1638.091
$ws.Range("L86").Value = 4089.889
$ws.Range("M86").Value = -515.0909999999999
$ws.Range("N86").Value = -6335.889
$ws.Range("H89").Value = 2741.4
$ws.Range("I89").Value = 1638.091
$ws.Range("J89").Value = 4089.889
$ws.Range("K89").Value = 8190.455
$ws.Range("L89").Value = 20449.445
$ws.Range("M89").Value = -2574.455
$ws.Range("N89").Value = -31681.445
$ws.Range("H107").Value = 1785.7
$ws.Range("I107").Value = 1528.5714
$ws.Range("J107").Value = 1924.1538
$ws.Range("K107").Value = 1528.5714
$ws.Range("L107").Value = 1924.1538
$ws.Range("M107").Value = 391.4286
$ws.Range("N107").Value = -5764.1538
$ws.Range("H134").Value = 2875016
$ws.Range("I134").Value = 1430.6666
$ws.Range("J134").Value = 41668416
$ws.Range("K134").Value = 4291.9998
$ws.Range("L134").Value = 125005248
$ws.Range("M134").Value = -1756.9998
$ws.Range("N134").Value = -125010318

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2906.9
$ws.Range("I132").Value = 2679.7058
$ws.Range("K132").Value = 8039.117400000001
$ws.Range("M132").Value = -5509.117400000001
$ws.Range("H134").Value = 3801.676
$ws.Range("I134").Value = 2621.8
$ws.Range("K134").Value = 7865.400000000001
$ws.Range("M134").Value = -5330.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 571
$ws.Range("I7").Value = 699.25
$ws.Range("K7").Value = 2097.75
$ws.Range("M7").Value = -1985.75
$ws.Range("H23").Value = 310.84616
$ws.Range("I23").Value = 338
$ws.Range("J23").Value = 293.875
$ws.Range("K23").Value = 1014
$ws.Range("L23").Value = 881.625
$ws.Range("M23").Value = -779
$ws.Range("N23").Value = -1351.625
$ws.Range("H39").Value = 20960.867
$ws.Range("J39").Value = 31151.4
$ws.Range("L39").Value = 93454.20000000001
$ws.Range("N39").Value = -94042.20000000001
$ws.Range("H55").Value = 143333700
$ws.Range("I55").Value = 210000540
$ws.Range("K55").Value = 630001620
$ws.Range("M55").Value = -630001443
$ws.Range("H57").Value = 4000
$ws.Range("I57").Value = 4000
$ws.Range("K57").Value = 12000
$ws.Range("M57").Value = -11441
$ws.Range("H87").Value = 16692.77
$ws.Range("I87").Value = 5601.2
$ws.Range("K87").Value = 16803.6
$ws.Range("M87").Value = -15555.6
$ws.Range("H90").Value = 16692.77
$ws.Range("I90").Value = 5601.2
$ws.Range("K90").Value = 50410.8
$ws.Range("M90").Value = -44170.8
$ws.Range("H131").Value = 4368.653
$ws.Range("J131").Value = 5729.778
$ws.Range("L131").Value = 17189.334
$ws.Range("N131").Value = -27269.334
$ws.Range("H132").Value = 1599.7142
$ws.Range("J132").Value = 1599.625
$ws.Range("L132").Value = 14396.625
$ws.Range("N132").Value = -19456.625
$ws.Range("H134").Value = 5342.136
$ws.Range("I134").Value = 1845.4375
$ws.Range("K134").Value = 5536.3125
$ws.Range("M134").Value = -466.3125
$ws.Range("H136").Value = 9537.333000000001
$ws.Range("I136").Value = 3343.3333
$ws.Range("J136").Value = 13666.667
$ws.Range("K136").Value = 10029.9999
$ws.Range("L136").Value = 41000.001
$ws.Range("M136").Value = -4929.999899999999
$ws.Range("N136").Value = -51200.001
$ws.Range("H137").Value = 9534
$ws.Range("I137").Value = 3971.6667
$ws.Range("J137").Value = 11202.7
$ws.Range("K137").Value = 11915.0001
$ws.Range("L137").Value = 33608.10000000001
$ws.Range("M137").Value = -6815.000100000001
$ws.Range("N137").Value = -43808.10000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 965.30554
$ws.Range("I97").Value = 759.62067
$ws.Range("J97").Value = 1817.4286
$ws.Range("K97").Value = 759.62067
$ws.Range("L97").Value = 1817.4286
$ws.Range("M97").Value = -263.62067
$ws.Range("N97").Value = -2809.4286
$ws.Range("H113").Value = 1357.5555
$ws.Range("I113").Value = 1357.5555
$ws.Range("K113").Value = 1357.5555
$ws.Range("M113").Value = 812.4445000000001
$ws.Range("H122").Value = 2656.423
$ws.Range("I122").Value = 2666.9583
$ws.Range("J122").Value = 2530
$ws.Range("K122").Value = 8000.874899999999
$ws.Range("L122").Value = 7590
$ws.Range("M122").Value = -5550.874899999999
$ws.Range("N122").Value = -12490
$ws.Range("H132").Value = 15677.96
$ws.Range("I132").Value = 8632.549999999999
$ws.Range("K132").Value = 25897.65
$ws.Range("M132").Value = -23367.65
$ws.Range("H136").Value = 76162.5
$ws.Range("J136").Value = 76162.5
$ws.Range("L136").Value = 228487.5
$ws.Range("N136").Value = -233587.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4299.45
$ws.Range("J46").Value = 5399.4
$ws.Range("L46").Value = 5399.4
$ws.Range("N46").Value = -5775.4
$ws.Range("H55").Value = 1474.091
$ws.Range("I55").Value = 1361.4667
$ws.Range("K55").Value = 1361.4667
$ws.Range("M55").Value = -1188.4667
$ws.Range("H82").Value = 3333.3635
$ws.Range("I82").Value = 3249.8333
$ws.Range("J82").Value = 3433.6
$ws.Range("K82").Value = 3249.8333
$ws.Range("L82").Value = 3433.6
$ws.Range("M82").Value = -2888.8333
$ws.Range("N82").Value = -4155.6
$ws.Range("H85").Value = 3333.3635
$ws.Range("I85").Value = 3249.8333
$ws.Range("J85").Value = 3433.6
$ws.Range("K85").Value = 3249.8333
$ws.Range("L85").Value = 3433.6
$ws.Range("M85").Value = -2001.8333
$ws.Range("N85").Value = -5929.6
$ws.Range("H93").Value = 3430.1
$ws.Range("I93").Value = 2680.9285
$ws.Range("K93").Value = 2680.9285
$ws.Range("M93").Value = -1432.9285
$ws.Range("H132").Value = 1518646.8
$ws.Range("I132").Value = 3033734
$ws.Range("K132").Value = 9101202
$ws.Range("M132").Value = -9098672
$ws.Range("H136").Value = 25003542
$ws.Range("I136").Value = 13892824
$ws.Range("K136").Value = 41678472
$ws.Range("M136").Value = -41675922

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 44996
$ws.Range("J45").Value = 44994.668
$ws.Range("L45").Value = 44994.668
$ws.Range("N45").Value = -45976.668
$ws.Range("H132").Value = 7939280
$ws.Range("I132").Value = 8774576
$ws.Range("J132").Value = 3971.5
$ws.Range("K132").Value = 26323728
$ws.Range("M132").Value = -26321198
$ws.Range("N132").Value = -16974.5
